$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect temporarily so values can be written
$ws.Unprotect()

# Update the confidential disclosure date from 2021-05-27 to 2021-05-28
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15
$ws.Range("D2").Value = 0.05687431591211282
$ws.Range("E2").Value = 0.001699430461142848

$ws.Range("D3").Value = 0.02389507561211988
$ws.Range("E3").Value = 0.002869714941649004

$ws.Range("D4").Value = 0.03150396150560723
$ws.Range("E4").Value = 0.005223880597014841

$ws.Range("D5").Value = 0.03234837738924633
$ws.Range("E5").Value = 0.003068661296509401

$ws.Range("D6").Value = 0.03707593081667359
$ws.Range("E6").Value = 0.005072923272035235

$ws.Range("D7").Value = 0.01877502880026332
$ws.Range("E7").Value = 0.00243486729973208

$ws.Range("D8").Value = 0.004468938818052866
$ws.Range("E8").Value = 0.004822446295484673

$ws.Range("D9").Value = 0.006935573615040398
$ws.Range("E9").Value = -0.0003766478342750457

$ws.Range("D10").Value = 0.07382271576687069
$ws.Range("E10").Value = 0.00424628450106157

$ws.Range("D11").Value = 0.07394026786204086
$ws.Range("E11").Value = 0.003179650238473775

$ws.Range("D12").Value = 0.1448816511628515
$ws.Range("E12").Value = -0.001550611229310217

$ws.Range("D13").Value = 0.3810588309052295
$ws.Range("E13").Value = -0.00008740494711989566

$ws.Range("D14").Value = 0.1144193318338911
$ws.Range("E14").Value = 0.004674577490111131

$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0.001507279086959912

# Restore sheet protection to match original state
$ws.Protect()
